$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture values that need to move before we overwrite/clear their source cells
$analisar = $ws.Range("D4").Value()
$darIdeias = $ws.Range("E4").Value()
$jogar = $ws.Range("E5").Value()
$pesquisar = $ws.Range("E6").Value()

# Clear the old "User Storys" column (B) entries - keep formatting, clear contents
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()

# Clear D4 (its value moved to E4)
$ws.Range("D4").ClearContents()

# E4 now holds what used to be in D4 ("Analisar o código dado")
$ws.Range("E4").Value = $analisar

# E5 and E6 are cleared (their values moved down into column G)
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()

# Append new rows to the "Done:" column (G)
$ws.Range("G9").Value = "Fazer pdf dos 3 User Stories mais votados pela equipa e submeter no moodle"
$ws.Range("G10").Value = $pesquisar
$ws.Range("G11").Value = $jogar
$ws.Range("G12").Value = $darIdeias

# Widen column G to fit the new, longer content
$ws.Range("G1").ColumnWidth = 68.02213541666667

# Update the active selection to reflect where the user ended up editing
$ws.Range("G15").Select()
